# Updates the crypto market snapshot on Sheet1 with refreshed
# prices, 1h volume changes, and a couple of re-ranked rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = "'98.385.90"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.15%  '

# Row 3: Ethereum
$ws.Range('D3').Value = "'3.419.25"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.35%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.02%  '

# Row 5: Solana
$ws.Range('D5').Value = "'256.38"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.54%  '

# Row 6: BNB
$ws.Range('D6').Value = "'668.52"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.39%  '

# Row 7: XRP
$ws.Range('D7').Value = "'1.47"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.00%  '

# Row 8: Dogecoin
$ws.Range('D8').Value = "'0.437"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.62%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  -2.43%  '

# Row 10: USDC
$ws.Range('E10').Value = '  -0.01%  '

# Row 11: LidoStakedEther
$ws.Range('D11').Value = "'3.416.92"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.41%  '

# Row 12: TRON
$ws.Range('E12').Value = '  +3.32%  '

# Row 13: Avalanche
$ws.Range('D13').Value = "'42.23"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.60%  '

# Row 14: Toncoin
$ws.Range('D14').Value = "'6.45"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +14.65%  '

# Row 15: WrappedBTC
$ws.Range('D15').Value = "'98.105.81"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '

# Row 16: ShibaInu
$ws.Range('D16').Value = "'0.0000268"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.74%  '

# Row 17: WrappedliquidstakedEther2.0
$ws.Range('D17').Value = "'4.051.14"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.14%  '

# Row 18: Polkadot
$ws.Range('D18').Value = "'9.08"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +20.80%  '

# Row 19: Stellar
$ws.Range('D19').Value = "'0.592"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +34.87%  '

# Row 20: WrappedEther
$ws.Range('D20').Value = "'3.403.94"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.89%  '

# Row 21: Chainlink
$ws.Range('D21').Value = "'17.72"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.93%  '

# Row 22: Uniswap
$ws.Range('D22').Value = "'11.06"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.91%  '

# Row 23: SuiNetwork
$ws.Range('D23').Value = "'3.46"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.00%  '

# Row 24: BitcoinCash
$ws.Range('D24').Value = "'512.43"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.61%  '

# Row 25: PEPE
$ws.Range('E25').Value = '  -1.33%  '

# Row 26: NEARProtocol
$ws.Range('E26').Value = '  +5.12%  '

# Row 27: Litecoin
$ws.Range('D27').Value = "'101.80"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.29%  '

# Row 28: Aptos
$ws.Range('D28').Value = "'12.87"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.44%  '

# Row 29: WrappedeETH
$ws.Range('D29').Value = "'3.599.30"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.27%  '

# Row 30: Hedera
$ws.Range('E30').Value = '  +1.48%  '

# Row 31: InternetComputer(DFINITY)
$ws.Range('D31').Value = "'11.63"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.20%  '

# Row 32: Cronos
$ws.Range('E32').Value = '  +2.95%  '

# Row 34: PancakeSwap
$ws.Range('D34').Value = "'2.50"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +17.69%  '

# Row 35: PolygonEcosystemToken
$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D35').Value = "'0.576"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.05%  '

# Row 36: Binance-PegBSC-USD
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').Value = "'0.999"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.16%  '

# Row 37: EthereumClassic
$ws.Range('D37').Value = "'30.11"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.88%  '

# Row 38: RenderToken
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = "'8.00"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.52%  '

# Row 39: Fetch.AI
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = "'1.51"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.87%  '

# Row 40: Bittensor
$ws.Range('D40').Value = "'538.92"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.24%  '

# Row 41: Kaspa
$ws.Range('D41').Value = "'0.155"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.97%  '

# Row 42: USDe
$ws.Range('E42').Value = '  -0.01%  '

# Row 43: ARBITRUM
$ws.Range('D43').Value = "'0.880"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.19%  '

# Row 44: WhiteBITCoin
$ws.Range('D44').Value = "'24.72"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.10%  '

# Row 45: MantraDAO
$ws.Range('D45').Value = "'3.82"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.53%  '

# Row 46: VeChain
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = "'0.0438"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.76%  '

# Row 47: Filecoin
$ws.Range('D47').Value = "'5.90"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.94%  '

# Row 48: Cosmos
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = "'8.99"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +13.55%  '

# Row 49: ImmutableX
$ws.Range('E49').Value = '  +14.14%  '

# Row 50: dogwifhat
$ws.Range('D50').Value = "'3.29"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.93%  '

# Row 51: OKB
$ws.Range('D51').Value = "'54.06"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.61%  '
